$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: shorten the VHSYS paragraph - drop trailing sentence
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute(
    " Mas nesse caso o sistema é online e o nosso cliente teria que fazer um plano mensal. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " ", 2)

Write-Output "change1 done"

# -----------------------------------------------------------------
# Change 3: merge the split run in the table cell and drop the
# stray bookmark that used to sit in the middle of "pedido pelo".
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Requisito técnico pedido pelo cliente",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Requisito técnico pedido pelo cliente", 2)

Write-Output "change3 done"

# -----------------------------------------------------------------
# Change 2: rewrite the "2.2. Justificativa" body paragraph
# -----------------------------------------------------------------

# locate the paragraph (the one right after the "2.2. Justificativa"
# heading) and left-align it
$headRng = $d.Content
$null = $headRng.Find.Execute("2.2. Justificativa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bodyPara = $headRng.Paragraphs(1).Next()
$bodyPara.Alignment = 0

# replace the whole sentence (in one shot) - the replacement inherits
# the language formatting ("pt-BR") that the matched run carried
$oldText = "Após o levantamento das soluções similares na subseção 2.1, aqui deve estar as comparações com as soluções  encontradas, justificando o porquê da solução proposta ser aderente com o cliente necessita."
$newText = "A solução encontrada no site da VHSYS é muito similar ao nosso projeto porém, eles utilizam o sistema ONLINE, onde o proprietário terá que ter o plano mensal, para ter acesso no site, além de que a empresa DracoTatto quer o sistema em LAN."

$bodyRange = $bodyPara.Range
$found2 = $bodyRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Output "change2-text found=$found2"

# bold the product name "VHSYS" (with its surrounding spaces)
$boldRng = $bodyPara.Range
$foundBold = $boldRng.Find.Execute(" VHSYS ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "change2-bold found=$foundBold"
if ($foundBold) {
    $boldRng.Font.Bold = 1
}

# drop the old "_GoBack" bookmark from the table cell (it moves to
# wherever the cursor ends up after the last text edit) and re-add it
# at the very end of the paragraph we just rewrote
$endRng = $bodyPara.Range
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRng)
Write-Output "change2-bookmark done"
